$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing cell formatting (style index 1) onto the new column C cells
$ws.Range("A1").Copy()
$ws.Range("C1:C3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add Russian translations in column C
$ws.Range("C1").Value = "человек"   # person
$ws.Range("C2").Value = "машина"    # car
$ws.Range("C3").Value = "грузовик"  # truck

# Remove the now-obsolete numeric values in column B (rows 2 and 3)
$ws.Range("B2").Clear()
$ws.Range("B3").Clear()
